{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: date + title line (contains a manual line break)\nparas.items[0].insertText(\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.10.24: \u26a1\ufe0f\ud83d\ude80\\u000bPredicting from Strings: Language Model Embeddings for Bayesian Optimization\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 1\nparas.items[1].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d4\u05e1\u05d5\u05d2 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05e1\u05d5\u05e7\u05e8 \u05d1\u05d3\u05e8\u05f4\u05db - \u05d0\u05d5\u05dc\u05d9 \u05de\u05ea\u05d5\u05da 300 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d9\u05e9 1-2 \u05db\u05d0\u05dc\u05d5 (\u05dc\u05d0 \u05d1\u05d8\u05d5\u05d7). \u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05e9\u05d4\u05e0\u05d5\u05e9\u05d0 \u05dc\u05d0 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05dc\u05d0 \u05e9\u05d9\u05e9 \u05e4\u05d7\u05d5\u05ea \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05d5 \u05d5\u05d4\u05d5\u05d0 \u05e0\u05d7\u05e9\u05d1 \u05e4\u05d7\u05d5\u05ea \u05f4\u05d1\u05d0\u05d6\u05d6\u05d9\u05f4 \u05dc\u05de\u05e8\u05d5\u05ea \u05d7\u05e9\u05d9\u05d1\u05d5\u05ea\u05d9. \u05db\u05de\u05d5 \u05e9\u05de\u05e9\u05ea\u05de\u05e2 \u05de\u05e9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05d5\u05e9\u05d0 \u05d4\u05d5\u05d0 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05d9\u05d9\u05e1\u05d9\u05d0\u05e0\u05d9\u05ea.\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 2\nparas.items[2].insertText(\n  \"\u05d1\u05d2\u05d3\u05d5\u05dc \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05d9\u05d9\u05e1\u05d9\u05d0\u05e0\u05d9\u05ea \u05d4\u05d9\u05d0 \u05d0\u05d7\u05d3 \u05d4\u05db\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d5\u05ea \u05ea\u05db\u05e0\u05d5\u05df \u05e0\u05d9\u05e1\u05d5\u05d9\u05dd \u05d5\u05dc\u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 black-box optimization \u05db\u05d0\u05e9\u05e8 \u05d4\u05d9\u05d0 \u05dc\u05de\u05e2\u05d6\u05e8 \u05d0\u05ea \u05de\u05d7\u05d9\u05e8 \u05e9\u05dc \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05d4\u05de\u05de\u05e7\u05e1\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4. \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d9\u05db\u05d5\u05dc\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05ea\u05e8\u05d5\u05e4\u05d4 (\u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea \u05d4\u05e8\u05db\u05d1\u05d4 \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9) \u05d0\u05d5 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05e8\u05e9\u05ea \u05d2\u05d3\u05d5\u05dc\u05d4. \u05d1\u05e9\u05e0\u05d9 \u05d4\u05de\u05e7\u05e8\u05d9\u05dd \u05db\u05dc \u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05e0\u05d4 \u05d9\u05e7\u05e8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d5\u05d9\u05e9 \u05e6\u05d5\u05e8\u05da \u05dc\u05de\u05d6\u05e2\u05e8 \u05d0\u05ea \u05db\u05de\u05d5\u05ea \u05d4\u05e4\u05e2\u05de\u05d9\u05dd \u05e9\u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05d5\u05ea\u05d4 (\u05dc\u05d1\u05d3\u05d9\u05e7\u05d4 \u05d4\u05e8\u05db\u05d1 \u05e9\u05dc \u05ea\u05e8\u05d5\u05e4\u05d4 \u05d0\u05d5 \u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05d4\u05d9\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05dd \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea).\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 3\nparas.items[3].insertText(\n  \"\u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05dc\u05d0 \u05de\u05e2\u05d8 \u05e9\u05d9\u05d8\u05d5\u05ea \u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d1\u05d7\u05d9\u05e8\u05ea \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea x \u05dc\u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05de\u05e6\u05d3 \u05d0\u05d7\u05d3 \u05d1\u05d5\u05d7\u05e8\u05ea \u05d0\u05d9\u05d6\u05d5\u05e8\u05d9\u05dd \u05d1\u05d4\u05dd \u05dc\u05d0 \u05d1\u05d3\u05e7\u05e0\u05d5 (exploration) \u05d5\u05de\u05e6\u05d3 \u05e9\u05e0\u05d9 \u05d2\u05dd \u05de\u05e0\u05e6\u05dc\u05ea \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2 \u05e9\u05dc\u05e0\u05d5 \u05e2\u05dc \u05e2\u05e8\u05db\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d1\u05d0\u05d9\u05d6\u05d5\u05e8\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d1\u05d9\u05e7\u05e8\u05e0\u05d5 (exploitation) \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05de\u05e6\u05d5\u05d0 \u05e0\u05e7\u05d5\u05d3\u05ea \u05de\u05e7\u05e1\u05d9\u05de\u05d5\u05dd \u05d8\u05d5\u05d1\u05d4 \u05d1\u05de\u05d0\u05de\u05e5 \u05de\u05d9\u05e0\u05d9\u05de\u05dc\u05d9. \u05e8\u05d5\u05d1 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05de\u05e0\u05e1\u05d5\u05ea \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 surrogate objective \u05d0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05de\u05d8\u05e8\u05d4 \u05d3\u05de\u05d4 \u05d4\u05d6\u05d5\u05dc\u05d4 \u05dc\u05d4\u05e4\u05e2\u05dc\u05d4 \u05db\u05d3\u05d9 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea x \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05e4\u05e2\u05dc\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05d6\u05d5\u05d2\u05d5\u05ea x \u05d5- (y=f(x)). \u05d4\u05d3\u05e8\u05da \u05d4\u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d9\u05d0 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05ea\u05d4\u05dc\u05d9\u05db\u05d9 \u05d2\u05d0\u05d5\u05e1 \u05db\u05d3\u05d9 \u05dc\u05de\u05d3\u05dc \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05de\u05d8\u05e8\u05d4 \u05d3\u05de\u05d4 \u05d5\u05d1\u05e2\u05d6\u05e8\u05ea\u05d4 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea \u05d4-x \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9.\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 4\nparas.items[4].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05e8\u05ea\u05d5\u05dd \u05d0\u05ea \u05d4-LLMs \u05dc\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4 \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05e9\u05e2\u05e8\u05da \u05d0\u05ea \u05d4\u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05d0\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc (f(x \u05e2\u05d1\u05d5\u05e8 x \u05e0\u05ea\u05d5\u05df. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05e4\u05db\u05d9\u05dd \u05d0\u05ea \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05e9\u05dc x \u05d5-y \u05d4\u05d9\u05d3\u05d5\u05e2\u05d9\u05dd \u05dc\u05e4\u05d5\u05e8\u05de\u05d8 \u05e9\u05dc string (\u05e0\u05d2\u05d9\u05d3 \u05dcjson \u05d4\u05de\u05db\u05d9\u05dc \u05d0\u05ea \u05e9\u05de\u05d5\u05ea \u05d4\u05e4\u05d9\u05e6'\u05e8\u05d9\u05dd \u05d5\u05d4\u05e2\u05e8\u05db\u05d9\u05dd \u05e9\u05dc\u05d4\u05dd). \u05dc\u05d0\u05d7\u05e8\u05d9 \u05de\u05db\u05df \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05de\u05d1\u05d5\u05e1\u05e1 LLMs \u05d4\u05de\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d0\u05dc\u05d5 \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d3\u05d9 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05e2\u05e8\u05da \u05e9\u05dc x \u05e9\u05e2\u05d1\u05d5\u05e8\u05d5 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d7\u05e9\u05d1 \u05d0\u05ea (f(x  (\u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05e9\u05d5\u05e0\u05d5\u05ea). \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 (\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d0 \u05de\u05d0\u05d5\u05de\u05df) \u05e2\u05dc \u05e1\u05d3\u05e8\u05d5\u05ea \u05f4\u05d6\u05d4\u05d1\u05f4 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea x \u05d5- (f(x \u05dc\u05de\u05e1\u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea. \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df k \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05e2\u05e8\u05da \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e2\u05d1\u05d5\u05e8 x_k+1 \u05dc k-\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd.\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 5\nparas.items[5].insertText(\n  \"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e0\u05d9\u05d7 \u05db\u05d9 \u05d0\u05ea \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05e2\u05e8\u05db\u05d9 \u05d4- x-\u05d9\u05dd \u05dc\u05d1\u05d3\u05d9\u05e7\u05d4 \u05de\u05ea\u05e7\u05d1\u05dc\u05d9\u05dd \u05d3\u05e8\u05da \u05d0\u05d9\u05d6\u05d4 \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d5\u05e0\u05d9 \u05e0\u05ea\u05d5\u05df.\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 6 (\"\u05d0\u05d4\u05d1\u05ea\u05d9 - \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05db\u05ea\u05d5\u05d1 \u05d9\u05e4\u05d4 \u05d5\u05d1\u05e8\u05d5\u05e8.\") is removed entirely\nparas.items[6].delete();\n\n// Paragraph 7 (the arxiv link) gets its URL updated\nparas.items[7].insertText(\n  \"https://arxiv.org/pdf/2410.10190\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1: date + title line (contains a manual line break, char 11)\n$d.Paragraphs(1).Range.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.10.24: \u26a1\ufe0f\ud83d\ude80\" + [char]11 + \"Predicting from Strings: Language Model Embeddings for Bayesian Optimization\"\n\n# Paragraph 2\n$d.Paragraphs(2).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d4\u05e1\u05d5\u05d2 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05e1\u05d5\u05e7\u05e8 \u05d1\u05d3\u05e8\u05f4\u05db - \u05d0\u05d5\u05dc\u05d9 \u05de\u05ea\u05d5\u05da 300 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d9\u05e9 1-2 \u05db\u05d0\u05dc\u05d5 (\u05dc\u05d0 \u05d1\u05d8\u05d5\u05d7). \u05dc\u05d0 \u05d1\u05d2\u05dc\u05dc \u05e9\u05d4\u05e0\u05d5\u05e9\u05d0 \u05dc\u05d0 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d0\u05dc\u05d0 \u05e9\u05d9\u05e9 \u05e4\u05d7\u05d5\u05ea \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05d5 \u05d5\u05d4\u05d5\u05d0 \u05e0\u05d7\u05e9\u05d1 \u05e4\u05d7\u05d5\u05ea \u05f4\u05d1\u05d0\u05d6\u05d6\u05d9\u05f4 \u05dc\u05de\u05e8\u05d5\u05ea \u05d7\u05e9\u05d9\u05d1\u05d5\u05ea\u05d9. \u05db\u05de\u05d5 \u05e9\u05de\u05e9\u05ea\u05de\u05e2 \u05de\u05e9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05d5\u05e9\u05d0 \u05d4\u05d5\u05d0 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05d9\u05d9\u05e1\u05d9\u05d0\u05e0\u05d9\u05ea.\"\n\n# Paragraph 3\n$d.Paragraphs(3).Range.Text = \"\u05d1\u05d2\u05d3\u05d5\u05dc \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d1\u05d9\u05d9\u05e1\u05d9\u05d0\u05e0\u05d9\u05ea \u05d4\u05d9\u05d0 \u05d0\u05d7\u05d3 \u05d4\u05db\u05dc\u05d9\u05dd \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d5\u05ea \u05ea\u05db\u05e0\u05d5\u05df \u05e0\u05d9\u05e1\u05d5\u05d9\u05dd \u05d5\u05dc\u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 black-box optimization \u05db\u05d0\u05e9\u05e8 \u05d4\u05d9\u05d0 \u05dc\u05de\u05e2\u05d6\u05e8 \u05d0\u05ea \u05de\u05d7\u05d9\u05e8 \u05e9\u05dc \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05d4\u05de\u05de\u05e7\u05e1\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4. \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d9\u05db\u05d5\u05dc\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05ea\u05e8\u05d5\u05e4\u05d4 (\u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea \u05d4\u05e8\u05db\u05d1\u05d4 \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9) \u05d0\u05d5 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05e8\u05e9\u05ea \u05d2\u05d3\u05d5\u05dc\u05d4. \u05d1\u05e9\u05e0\u05d9 \u05d4\u05de\u05e7\u05e8\u05d9\u05dd \u05db\u05dc \u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05e0\u05d4 \u05d9\u05e7\u05e8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d5\u05d9\u05e9 \u05e6\u05d5\u05e8\u05da \u05dc\u05de\u05d6\u05e2\u05e8 \u05d0\u05ea \u05db\u05de\u05d5\u05ea \u05d4\u05e4\u05e2\u05de\u05d9\u05dd \u05e9\u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05d5\u05ea\u05d4 (\u05dc\u05d1\u05d3\u05d9\u05e7\u05d4 \u05d4\u05e8\u05db\u05d1 \u05e9\u05dc \u05ea\u05e8\u05d5\u05e4\u05d4 \u05d0\u05d5 \u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05e9\u05d9\u05dc\u05d5\u05d1 \u05d4\u05d9\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05dd \u05e9\u05dc \u05d4\u05e8\u05e9\u05ea).\"\n\n# Paragraph 4\n$d.Paragraphs(4).Range.Text = \"\u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05dc\u05d0 \u05de\u05e2\u05d8 \u05e9\u05d9\u05d8\u05d5\u05ea \u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d1\u05d7\u05d9\u05e8\u05ea \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea x \u05dc\u05d0\u05d1\u05dc\u05d5\u05d0\u05e6\u05d9\u05d4 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05de\u05e6\u05d3 \u05d0\u05d7\u05d3 \u05d1\u05d5\u05d7\u05e8\u05ea \u05d0\u05d9\u05d6\u05d5\u05e8\u05d9\u05dd \u05d1\u05d4\u05dd \u05dc\u05d0 \u05d1\u05d3\u05e7\u05e0\u05d5 (exploration) \u05d5\u05de\u05e6\u05d3 \u05e9\u05e0\u05d9 \u05d2\u05dd \u05de\u05e0\u05e6\u05dc\u05ea \u05d0\u05ea \u05d4\u05d9\u05d3\u05e2 \u05e9\u05dc\u05e0\u05d5 \u05e2\u05dc \u05e2\u05e8\u05db\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05de\u05d8\u05e8\u05d4 \u05d1\u05d0\u05d9\u05d6\u05d5\u05e8\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d1\u05d9\u05e7\u05e8\u05e0\u05d5 (exploitation) \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05de\u05e6\u05d5\u05d0 \u05e0\u05e7\u05d5\u05d3\u05ea \u05de\u05e7\u05e1\u05d9\u05de\u05d5\u05dd \u05d8\u05d5\u05d1\u05d4 \u05d1\u05de\u05d0\u05de\u05e5 \u05de\u05d9\u05e0\u05d9\u05de\u05dc\u05d9. \u05e8\u05d5\u05d1 \u05d4\u05e9\u05d9\u05d8\u05d5\u05ea \u05de\u05e0\u05e1\u05d5\u05ea \u05dc\u05d1\u05e0\u05d5\u05ea \u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 surrogate objective \u05d0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05de\u05d8\u05e8\u05d4 \u05d3\u05de\u05d4 \u05d4\u05d6\u05d5\u05dc\u05d4 \u05dc\u05d4\u05e4\u05e2\u05dc\u05d4 \u05db\u05d3\u05d9 \u05dc\u05de\u05e6\u05d5\u05d0 \u05d0\u05ea x \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05e4\u05e2\u05dc\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea (\u05db\u05dc\u05d5\u05de\u05e8 \u05d6\u05d5\u05d2\u05d5\u05ea x \u05d5- (y=f(x)). \u05d4\u05d3\u05e8\u05da \u05d4\u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d9\u05d0 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05ea\u05d4\u05dc\u05d9\u05db\u05d9 \u05d2\u05d0\u05d5\u05e1 \u05db\u05d3\u05d9 \u05dc\u05de\u05d3\u05dc \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05de\u05d8\u05e8\u05d4 \u05d3\u05de\u05d4 \u05d5\u05d1\u05e2\u05d6\u05e8\u05ea\u05d4 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea \u05d4-x \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9.\"\n\n# Paragraph 5\n$d.Paragraphs(5).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05e8\u05ea\u05d5\u05dd \u05d0\u05ea \u05d4-LLMs \u05dc\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4 \u05d1\u05de\u05d8\u05e8\u05d4 \u05dc\u05e9\u05e2\u05e8\u05da \u05d0\u05ea \u05d4\u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05d0\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc (f(x \u05e2\u05d1\u05d5\u05e8 x \u05e0\u05ea\u05d5\u05df. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05e4\u05db\u05d9\u05dd \u05d0\u05ea \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05e9\u05dc x \u05d5-y \u05d4\u05d9\u05d3\u05d5\u05e2\u05d9\u05dd \u05dc\u05e4\u05d5\u05e8\u05de\u05d8 \u05e9\u05dc string (\u05e0\u05d2\u05d9\u05d3 \u05dcjson \u05d4\u05de\u05db\u05d9\u05dc \u05d0\u05ea \u05e9\u05de\u05d5\u05ea \u05d4\u05e4\u05d9\u05e6'\u05e8\u05d9\u05dd \u05d5\u05d4\u05e2\u05e8\u05db\u05d9\u05dd \u05e9\u05dc\u05d4\u05dd). \u05dc\u05d0\u05d7\u05e8\u05d9 \u05de\u05db\u05df \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05de\u05d1\u05d5\u05e1\u05e1 LLMs \u05d4\u05de\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05de\u05db\u05e0\u05d9\u05e1\u05d9\u05dd \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d0\u05dc\u05d5 \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d3\u05d9 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05e2\u05e8\u05da \u05e9\u05dc x \u05e9\u05e2\u05d1\u05d5\u05e8\u05d5 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d7\u05e9\u05d1 \u05d0\u05ea (f(x  (\u05ea\u05d5\u05d7\u05dc\u05ea \u05d5\u05e9\u05d5\u05e0\u05d5\u05ea). \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 (\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d0 \u05de\u05d0\u05d5\u05de\u05df) \u05e2\u05dc \u05e1\u05d3\u05e8\u05d5\u05ea \u05f4\u05d6\u05d4\u05d1\u05f4 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea x \u05d5- (f(x \u05dc\u05de\u05e1\u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea. \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df k \u05d4\u05d6\u05d5\u05d2\u05d5\u05ea \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05e2\u05e8\u05da \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e2\u05d1\u05d5\u05e8 x_k+1 \u05dc k-\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd.\"\n\n# Paragraph 6\n$d.Paragraphs(6).Range.Text = \"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e0\u05d9\u05d7 \u05db\u05d9 \u05d0\u05ea \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05e2\u05e8\u05db\u05d9 \u05d4- x-\u05d9\u05dd \u05dc\u05d1\u05d3\u05d9\u05e7\u05d4 \u05de\u05ea\u05e7\u05d1\u05dc\u05d9\u05dd \u05d3\u05e8\u05da \u05d0\u05d9\u05d6\u05d4 \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d5\u05e0\u05d9 \u05e0\u05ea\u05d5\u05df.\"\n\n# Paragraph 7 (\"\u05d0\u05d4\u05d1\u05ea\u05d9 - \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05db\u05ea\u05d5\u05d1 \u05d9\u05e4\u05d4 \u05d5\u05d1\u05e8\u05d5\u05e8.\") is removed entirely\n$d.Paragraphs(7).Range.Delete()\n\n# Paragraph 8 (now the last paragraph, the arxiv link) gets its URL updated\n$d.Paragraphs(7).Range.Text = \"https://arxiv.org/pdf/2410.10190\"\n"}
